# edit.ps1
# Applies the "update scripts with new tpm" change to the Ntf5-Ntrk2 NATMI
# ligand-receptor results sheet:
#   - A new sending/target cluster "MuSCs" is added, extending the
#     sending-cluster x target-cluster grid from 2x3 (6 rows) to 3x4 (12 rows).
#   - All numeric (TPM-derived) statistics are refreshed with newly computed
#     values across every existing and new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("A2").Value = "FAPs"
    $ws.Range("B2").Value = "Ntf5"
    $ws.Range("C2").Value = "Ntrk2"
    $ws.Range("D2").Value = "ECs"
    $ws.Range("E2").Value = 1
    $ws.Range("F2").Value = 0.3333333333333333
    $ws.Range("G2").Value = 0.563831
    $ws.Range("H2").Value = 1.691493
    $ws.Range("I2").Value = 0.5426742997499485
    $ws.Range("J2").Value = 0.5426742997499486
    $ws.Range("K2").Value = 3
    $ws.Range("L2").Value = 1
    $ws.Range("M2").Value = 0.8630909999999999
    $ws.Range("N2").Value = 2.589273
    $ws.Range("O2").Value = 0.029079428547613
    $ws.Range("P2").Value = 0.029079428547613
    $ws.Range("Q2").Value = 0.486637461621
    $ws.Range("R2").Value = 4.379737154589
    $ws.Range("S2").Value = 0.01578065852420455
    $ws.Range("T2").Value = 0.01578065852420455

    # Row 3
    $ws.Range("A3").Value = "FAPs"
    $ws.Range("B3").Value = "Ntf5"
    $ws.Range("C3").Value = "Ntrk2"
    $ws.Range("D3").Value = "FAPs"
    $ws.Range("E3").Value = 1
    $ws.Range("F3").Value = 0.3333333333333333
    $ws.Range("G3").Value = 0.563831
    $ws.Range("H3").Value = 1.691493
    $ws.Range("I3").Value = 0.5426742997499485
    $ws.Range("J3").Value = 0.5426742997499486
    $ws.Range("K3").Value = 3
    $ws.Range("L3").Value = 1
    $ws.Range("M3").Value = 24.72809866666667
    $ws.Range("N3").Value = 74.184296
    $ws.Range("O3").Value = 0.833143872773158
    $ws.Range("P3").Value = 0.8331438727731579
    $ws.Range("Q3").Value = 13.94246859932533
    $ws.Range("R3").Value = 125.482217393928
    $ws.Range("S3").Value = 0.4521257677481337
    $ws.Range("T3").Value = 0.4521257677481338

    # Row 4
    $ws.Range("A4").Value = "FAPs"
    $ws.Range("B4").Value = "Ntf5"
    $ws.Range("C4").Value = "Ntrk2"
    $ws.Range("D4").Value = "MuSCs"
    $ws.Range("E4").Value = 1
    $ws.Range("F4").Value = 0.3333333333333333
    $ws.Range("G4").Value = 0.563831
    $ws.Range("H4").Value = 1.691493
    $ws.Range("I4").Value = 0.5426742997499485
    $ws.Range("J4").Value = 0.5426742997499486
    $ws.Range("K4").Value = 3
    $ws.Range("L4").Value = 1
    $ws.Range("M4").Value = 4.080109666666666
    $ws.Range("N4").Value = 12.240329
    $ws.Range("O4").Value = 0.137467842346008
    $ws.Range("P4").Value = 0.137467842346008
    $ws.Range("Q4").Value = 2.300492313466333
    $ws.Range("R4").Value = 20.704430821197
    $ws.Range("S4").Value = 0.07460026508325623
    $ws.Range("T4").Value = 0.07460026508325625

    # Row 5
    $ws.Range("A5").Value = "FAPs"
    $ws.Range("B5").Value = "Ntf5"
    $ws.Range("C5").Value = "Ntrk2"
    $ws.Range("D5").Value = "Resolving-Mac"
    $ws.Range("E5").Value = 1
    $ws.Range("F5").Value = 0.3333333333333333
    $ws.Range("G5").Value = 0.563831
    $ws.Range("H5").Value = 1.691493
    $ws.Range("I5").Value = 0.5426742997499485
    $ws.Range("J5").Value = 0.5426742997499486
    $ws.Range("K5").Value = 1
    $ws.Range("L5").Value = 0.3333333333333333
    $ws.Range("M5").Value = 0.009167
    $ws.Range("N5").Value = 0.027501
    $ws.Range("O5").Value = 0.0003088563332209099
    $ws.Range("P5").Value = 0.0003088563332209099
    $ws.Range("Q5").Value = 0.005168638777
    $ws.Range("R5").Value = 0.046517748993
    $ws.Range("S5").Value = 0.0001676083943539941
    $ws.Range("T5").Value = 0.0001676083943539941

    # Row 6
    $ws.Range("A6").Value = "MuSCs"
    $ws.Range("B6").Value = "Ntf5"
    $ws.Range("C6").Value = "Ntrk2"
    $ws.Range("D6").Value = "ECs"
    $ws.Range("E6").Value = 1
    $ws.Range("F6").Value = 0.3333333333333333
    $ws.Range("G6").Value = 0.05363299999999999
    $ws.Range("H6").Value = 0.160899
    $ws.Range("I6").Value = 0.05162052231695133
    $ws.Range("J6").Value = 0.05162052231695134
    $ws.Range("K6").Value = 3
    $ws.Range("L6").Value = 1
    $ws.Range("M6").Value = 0.8630909999999999
    $ws.Range("N6").Value = 2.589273
    $ws.Range("O6").Value = 0.029079428547613
    $ws.Range("P6").Value = 0.029079428547613
    $ws.Range("Q6").Value = 0.04629015960299999
    $ws.Range("R6").Value = 0.4166114364269999
    $ws.Range("S6").Value = 0.001501095290306248
    $ws.Range("T6").Value = 0.001501095290306249

    # Row 7
    $ws.Range("A7").Value = "MuSCs"
    $ws.Range("B7").Value = "Ntf5"
    $ws.Range("C7").Value = "Ntrk2"
    $ws.Range("D7").Value = "FAPs"
    $ws.Range("E7").Value = 1
    $ws.Range("F7").Value = 0.3333333333333333
    $ws.Range("G7").Value = 0.05363299999999999
    $ws.Range("H7").Value = 0.160899
    $ws.Range("I7").Value = 0.05162052231695133
    $ws.Range("J7").Value = 0.05162052231695134
    $ws.Range("K7").Value = 3
    $ws.Range("L7").Value = 1
    $ws.Range("M7").Value = 24.72809866666667
    $ws.Range("N7").Value = 74.184296
    $ws.Range("O7").Value = 0.833143872773158
    $ws.Range("P7").Value = 0.8331438727731579
    $ws.Range("Q7").Value = 1.326242115789333
    $ws.Range("R7").Value = 11.936179042104
    $ws.Range("S7").Value = 0.04300732187771807
    $ws.Range("T7").Value = 0.04300732187771807

    # Row 8
    $ws.Range("A8").Value = "MuSCs"
    $ws.Range("B8").Value = "Ntf5"
    $ws.Range("C8").Value = "Ntrk2"
    $ws.Range("D8").Value = "MuSCs"
    $ws.Range("E8").Value = 1
    $ws.Range("F8").Value = 0.3333333333333333
    $ws.Range("G8").Value = 0.05363299999999999
    $ws.Range("H8").Value = 0.160899
    $ws.Range("I8").Value = 0.05162052231695133
    $ws.Range("J8").Value = 0.05162052231695134
    $ws.Range("K8").Value = 3
    $ws.Range("L8").Value = 1
    $ws.Range("M8").Value = 4.080109666666666
    $ws.Range("N8").Value = 12.240329
    $ws.Range("O8").Value = 0.137467842346008
    $ws.Range("P8").Value = 0.137467842346008
    $ws.Range("Q8").Value = 0.2188285217523333
    $ws.Range("R8").Value = 1.969456695771
    $ws.Range("S8").Value = 0.007096161823685255
    $ws.Range("T8").Value = 0.007096161823685256

    # Row 9
    $ws.Range("A9").Value = "MuSCs"
    $ws.Range("B9").Value = "Ntf5"
    $ws.Range("C9").Value = "Ntrk2"
    $ws.Range("D9").Value = "Resolving-Mac"
    $ws.Range("E9").Value = 1
    $ws.Range("F9").Value = 0.3333333333333333
    $ws.Range("G9").Value = 0.05363299999999999
    $ws.Range("H9").Value = 0.160899
    $ws.Range("I9").Value = 0.05162052231695133
    $ws.Range("J9").Value = 0.05162052231695134
    $ws.Range("K9").Value = 1
    $ws.Range("L9").Value = 0.3333333333333333
    $ws.Range("M9").Value = 0.009167
    $ws.Range("N9").Value = 0.027501
    $ws.Range("O9").Value = 0.0003088563332209099
    $ws.Range("P9").Value = 0.0003088563332209099
    $ws.Range("Q9").Value = 0.000491653711
    $ws.Range("R9").Value = 0.004424883399
    $ws.Range("S9").Value = [double]"1.594332524176174E-05"
    $ws.Range("T9").Value = [double]"1.594332524176174E-05"

    # Row 10
    $ws.Range("A10").Value = "Resolving-Mac"
    $ws.Range("B10").Value = "Ntf5"
    $ws.Range("C10").Value = "Ntrk2"
    $ws.Range("D10").Value = "ECs"
    $ws.Range("E10").Value = 2
    $ws.Range("F10").Value = 0.6666666666666666
    $ws.Range("G10").Value = 0.421522
    $ws.Range("H10").Value = 1.264566
    $ws.Range("I10").Value = 0.4057051779331001
    $ws.Range("J10").Value = 0.4057051779331002
    $ws.Range("K10").Value = 3
    $ws.Range("L10").Value = 1
    $ws.Range("M10").Value = 0.8630909999999999
    $ws.Range("N10").Value = 2.589273
    $ws.Range("O10").Value = 0.029079428547613
    $ws.Range("P10").Value = 0.029079428547613
    $ws.Range("Q10").Value = 0.3638118445019999
    $ws.Range("R10").Value = 3.274306600517999
    $ws.Range("S10").Value = 0.0117976747331022
    $ws.Range("T10").Value = 0.0117976747331022

    # Row 11
    $ws.Range("A11").Value = "Resolving-Mac"
    $ws.Range("B11").Value = "Ntf5"
    $ws.Range("C11").Value = "Ntrk2"
    $ws.Range("D11").Value = "FAPs"
    $ws.Range("E11").Value = 2
    $ws.Range("F11").Value = 0.6666666666666666
    $ws.Range("G11").Value = 0.421522
    $ws.Range("H11").Value = 1.264566
    $ws.Range("I11").Value = 0.4057051779331001
    $ws.Range("J11").Value = 0.4057051779331002
    $ws.Range("K11").Value = 3
    $ws.Range("L11").Value = 1
    $ws.Range("M11").Value = 24.72809866666667
    $ws.Range("N11").Value = 74.184296
    $ws.Range("O11").Value = 0.833143872773158
    $ws.Range("P11").Value = 0.8331438727731579
    $ws.Range("Q11").Value = 10.42343760617067
    $ws.Range("R11").Value = 93.81093845553599
    $ws.Range("S11").Value = 0.3380107831473062
    $ws.Range("T11").Value = 0.3380107831473062

    # Row 12
    $ws.Range("A12").Value = "Resolving-Mac"
    $ws.Range("B12").Value = "Ntf5"
    $ws.Range("C12").Value = "Ntrk2"
    $ws.Range("D12").Value = "MuSCs"
    $ws.Range("E12").Value = 2
    $ws.Range("F12").Value = 0.6666666666666666
    $ws.Range("G12").Value = 0.421522
    $ws.Range("H12").Value = 1.264566
    $ws.Range("I12").Value = 0.4057051779331001
    $ws.Range("J12").Value = 0.4057051779331002
    $ws.Range("K12").Value = 3
    $ws.Range("L12").Value = 1
    $ws.Range("M12").Value = 4.080109666666666
    $ws.Range("N12").Value = 12.240329
    $ws.Range("O12").Value = 0.137467842346008
    $ws.Range("P12").Value = 0.137467842346008
    $ws.Range("Q12").Value = 1.719855986912666
    $ws.Range("R12").Value = 15.478703882214
    $ws.Range("S12").Value = 0.05577141543906654
    $ws.Range("T12").Value = 0.05577141543906655

    # Row 13
    $ws.Range("A13").Value = "Resolving-Mac"
    $ws.Range("B13").Value = "Ntf5"
    $ws.Range("C13").Value = "Ntrk2"
    $ws.Range("D13").Value = "Resolving-Mac"
    $ws.Range("E13").Value = 2
    $ws.Range("F13").Value = 0.6666666666666666
    $ws.Range("G13").Value = 0.421522
    $ws.Range("H13").Value = 1.264566
    $ws.Range("I13").Value = 0.4057051779331001
    $ws.Range("J13").Value = 0.4057051779331002
    $ws.Range("K13").Value = 1
    $ws.Range("L13").Value = 0.3333333333333333
    $ws.Range("M13").Value = 0.009167
    $ws.Range("N13").Value = 0.027501
    $ws.Range("O13").Value = 0.0003088563332209099
    $ws.Range("P13").Value = 0.0003088563332209099
    $ws.Range("Q13").Value = 0.003864092173999999
    $ws.Range("R13").Value = 0.034776829566
    $ws.Range("S13").Value = 0.0001253046136251541
    $ws.Range("T13").Value = 0.0001253046136251542

# Refresh the used range so the dimension reflects the newly written rows.
$null = $ws.UsedRange
